# Add file upload support with logging and API integration
# - appends a new "start" record to the "starts" log sheet
# - adds an "uploads" sheet logging the uploaded file
# - adds an "operations" sheet logging the start + upload operations

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "starts" sheet: append the new session's start record as row 4
# ---------------------------------------------------------------------------
$starts = $wb.Worksheets.Item("starts")
$starts.Range("A4").Value = "2025-08-13T14:45:15.474683"
$starts.Range("B4").Value = "2025-08-13T14:45:15.474654"
$starts.Range("C4").Value = "a8c1ea65-a4e7-456d-b788-2d96dc55cb51"
$starts.Range("D4").Value = "E777"

# ---------------------------------------------------------------------------
# 2. New "uploads" sheet (appended after the last existing sheet)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$uploads = $wb.Worksheets.Add($null, $lastSheet)
$uploads.Name = "uploads"

$uploads.Range("A1").Value = "write_time"
$uploads.Range("B1").Value = "timestamp"
$uploads.Range("C1").Value = "session_id"
$uploads.Range("D1").Value = "employee_id"
$uploads.Range("E1").Value = "step_number"
$uploads.Range("F1").Value = "field_name"
$uploads.Range("G1").Value = "original_filename"
$uploads.Range("H1").Value = "saved_path"
$uploads.Range("I1").Value = "file_size_bytes"
$uploads.Range("J1").Value = "mime_type"

$uploads.Range("A2").Value = "2025-08-13T14:45:15.535796"
$uploads.Range("B2").Value = "2025-08-13T14:45:15.535782"
$uploads.Range("C2").Value = "a8c1ea65-a4e7-456d-b788-2d96dc55cb51"
$uploads.Range("D2").Value = "E777"
$uploads.Range("E2").Value = 2
$uploads.Range("F2").Value = "file_hukou"
$uploads.Range("G2").Value = "hukou.png"
$uploads.Range("H2").Value = "logs/uploads/E777/a8c1ea65-a4e7-456d-b788-2d96dc55cb51/20250813T144515_hukou.png"
$uploads.Range("I2").Value = 5
$uploads.Range("J2").Value = "image/png"

# ---------------------------------------------------------------------------
# 3. New "operations" sheet (appended after "uploads")
# ---------------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$operations = $wb.Worksheets.Add($null, $lastSheet2)
$operations.Name = "operations"

$operations.Range("A1").Value = "write_time"
$operations.Range("B1").Value = "timestamp"
$operations.Range("C1").Value = "session_id"
$operations.Range("D1").Value = "employee_id"
$operations.Range("E1").Value = "operation_type"
$operations.Range("F1").Value = "step_number"
$operations.Range("G1").Value = "name"
$operations.Range("H1").Value = "value"
$operations.Range("I1").Value = "extra"

$operations.Range("A2").Value = "2025-08-13T14:45:15.503579"
$operations.Range("B2").Value = "2025-08-13T14:45:15.474654"
$operations.Range("C2").Value = "a8c1ea65-a4e7-456d-b788-2d96dc55cb51"
$operations.Range("D2").Value = "E777"
$operations.Range("E2").Value = "start"
$operations.Range("F2").Value = 0
$operations.Range("G2").Value = "start"

$operations.Range("A3").Value = "2025-08-13T14:45:15.562053"
$operations.Range("B3").Value = "2025-08-13T14:45:15.535782"
$operations.Range("C3").Value = "a8c1ea65-a4e7-456d-b788-2d96dc55cb51"
$operations.Range("D3").Value = "E777"
$operations.Range("E3").Value = "upload"
$operations.Range("F3").Value = 2
$operations.Range("G3").Value = "file_hukou"
$operations.Range("H3").Value = "hukou.png"
$operations.Range("I3").Value = "logs/uploads/E777/a8c1ea65-a4e7-456d-b788-2d96dc55cb51/20250813T144515_hukou.png"

# ---------------------------------------------------------------------------
# Restore the original active sheet/selection ("starts" stays active, as in
# the source workbook - only the sheet list and "starts" data changed).
# ---------------------------------------------------------------------------
$starts.Activate()
$null = $starts.Range("A1").Select()

